$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 09:54"

# --- Estados Unidos (row 4): new confirmed cases + recovered updated ---
$ws.Range("B4").Value = 4170131
$ws.Range("C4").Value = 140
$ws.Range("E4").Value = 2043181

# --- India (row 6) ---
$ws.Range("B6").Value = 1290284
$ws.Range("C6").Value = 2154
$ws.Range("D6").Value = 817738
$ws.Range("E6").Value = 441889
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 30657

# --- Rusia (row 7) ---
$ws.Range("B7").Value = 800849
$ws.Range("C7").Value = 5811
$ws.Range("D7").Value = 588774
$ws.Range("E7").Value = 199029
$ws.Range("G7").Value = 154
$ws.Range("H7").Value = 13046

# --- Kuwait / Ucrania (rows 38-39) swap ranking order and refresh figures ---
# Ucrania overtakes Kuwait in total cases, so it now occupies row 38.
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 62823
$ws.Range("C38").Value = 972
$ws.Range("D38").Value = 34886
$ws.Range("E38").Value = 26366
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 20
$ws.Range("H38").Value = 1571

$ws.Range("A39").Value = "Kuwait"
$ws.Range("B39").Value = 61872
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 52247
$ws.Range("E39").Value = 9204
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 421

# --- Singapur (row 46) ---
$ws.Range("B46").Value = 49375
$ws.Range("C46").Value = 277
$ws.Range("E46").Value = 4333

# --- Armenia (row 53) ---
$ws.Range("B53").Value = 36613
$ws.Range("C53").Value = 451
$ws.Range("D53").Value = 25734
$ws.Range("E53").Value = 10187
$ws.Range("G53").Value = 4
$ws.Range("H53").Value = 692

# --- Afganistan (row 54) ---
$ws.Range("B54").Value = 35981
$ws.Range("C54").Value = 53
$ws.Range("D54").Value = 24602
$ws.Range("E54").Value = 10154
$ws.Range("G54").Value = 14
$ws.Range("H54").Value = 1225

# --- Sri Lanka (row 114) ---
$ws.Range("D114").Value = 2094
$ws.Range("E114").Value = 648

# --- Estonia (row 125) ---
$ws.Range("B125").Value = 2028
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 1915
$ws.Range("E125").Value = 44

# --- Lituania (row 127) ---
$ws.Range("B127").Value = 1986
$ws.Range("C127").Value = 26
$ws.Range("D127").Value = 1616
$ws.Range("E127").Value = 290

# --- Letonia (row 139) ---
$ws.Range("B139").Value = 1205
$ws.Range("C139").Value = 2
$ws.Range("E139").Value = 129

# --- Islas Malvinas / Groenlandia (rows 210-211) swap display order (tie on totals) ---
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"
